# edit.ps1
# Updates the pl_mw.xlsx "Case_1_6" results sheet with new simulation values
# (commit: "case with 380 kV done").
#
# The sheet holds a results table in rows 2-25 (indices 0-23, column A),
# with result columns B, D, E, F, G, H, I, N updated to new simulation
# output values. Columns C, J, K, L, M, O remain 0 and are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B" = @(0.2159311196958242, 0.1913949238933697, 0.1762869034279646, 0.1701201473110672, 0.1690955695420513, 0.1762037765289648, 0.2074802658006831, 0.2684501888307409, 0.3129948669848943, 0.3331993617551348, 0.3408412057246721, 0.3391958157675106, 0.333828248206629, 0.3305392456624361, 0.3116732098856687, 0.3000838993120851, 0.2934125106845613, 0.2911527646054139, 0.3013181776852321, 0.3354050867557419, 0.3576291578720543, 0.3457728835230682, 0.3007601869815915, 0.2519980490933165)
    "D" = @(0.04416460176246062, 0.03830832642090343, 0.03472637085036467, 0.03326964694228707, 0.03302792226286044, 0.03470671369345268, 0.04214216942487781, 0.05685744970638495, 0.06778656051733378, 0.07279217337988086, 0.07469321267035411, 0.0742835337146488, 0.0729484585700817, 0.07213142522000737, 0.06746016902452823, 0.06460364233217319, 0.06296379515380579, 0.06240909297251562, 0.06490739369498044, 0.07334044737054057, 0.07888443913425647, 0.07592231310724173, 0.06477006025150445, 0.05285854641996934)
    "E" = @(0.4396191691992186, 0.3823479279035951, 0.3474053453013823, 0.333216280922386, 0.3308630769826664, 0.3472137901143384, 0.4198225367166089, 0.5642327170640442, 0.6719446234360191, 0.7213816670721798, 0.7401723167001393, 0.7361221926845474, 0.7229261437163217, 0.7148524956228215, 0.6687232001523142, 0.6405414693056741, 0.6243729184997164, 0.6189053326469747, 0.6435371872360207, 0.726800189790751, 0.7816282444840539, 0.752325552932561, 0.6421827193271525, 0.524910646920091)
    "F" = @(1.528290297222981, 1.363286054659909, 1.262945067803372, 1.222288067622259, 1.215550724464379, 1.262395826126067, 1.47118785694326, 1.888902018709359, 2.201659492699889, 2.345399064368451, 2.400055170825851, 2.388273756458403, 2.349891060477688, 2.326410295331868, 2.192296469740057, 2.110407240501445, 2.063443689084153, 2.047565729167701, 2.119110209985365, 2.361158763589458, 2.520670277776958, 2.435410421870188, 2.115175238321626, 1.774930254878313)
    "G" = @(1.673767081463382, 1.472460773835678, 1.349747030403535, 1.299949458491369, 1.291692836908965, 1.349074614279345, 1.604163365212003, 2.112127782647462, 2.491042104895939, 2.664885521213535, 2.730945692543116, 2.71670796254466, 2.670315630292464, 2.641929398215893, 2.479712166869206, 2.38058685249996, 2.323710192705676, 2.304475839262807, 2.391124544869683, 2.683935813114317, 2.876650817821144, 2.773666093207737, 2.386360107228029, 1.973780395381823)
    "H" = @(1.010542361786406, 0.9231199012414208, 0.8702008080870769, 0.8488200623750686, 0.8452807431618226, 0.8699117225502562, 0.98023762357343, 1.202918706961327, 1.370841485233257, 1.44827547428082, 1.477756451588164, 1.471400018306724, 1.450697669861711, 1.438037758433723, 1.365802710240757, 1.321762217021217, 1.2965291561203, 1.288002240429705, 1.326440215915113, 1.45677408627796, 1.542882076632281, 1.496837023020134, 1.324325024093469, 1.14195305260256)
    "I" = @(0.3757315010515043, 0.3656238581545637, 0.3596127975951688, 0.3572128864356259, 0.3568173988916286, 0.359580229755295, 0.3722062484302171, 0.3984893586954925, 0.4186982167615625, 0.4280806917633271, 0.431660259156402, 0.4308881605543888, 0.4283746547068858, 0.4268385096919047, 0.4180888088535966, 0.4127692216154415, 0.4097273953840102, 0.4087005711172083, 0.4133336567285752, 0.429112214418673, 0.4395793436413413, 0.4339788619097504, 0.4130784242511822, 0.3912191043998661)
    "N" = @(1.538965158849294, 1.440536212657179, 1.380350891855613, 1.3558914054328, 1.351834048475041, 1.380020748618335, 1.504976964631908, 1.75186015072677, 1.93419659562025, 2.017317443662478, 2.048815004509208, 2.042030543333738, 2.01990835705044, 2.006360570313319, 1.928767661017844, 1.881209144293734, 1.85387144156806, 1.844618345498816, 1.88627012809178, 2.026405625114307, 2.118115898853603, 2.069158297353852, 1.883982043575145, 1.684892220676289)
}

foreach ($col in $data.Keys) {
    $values = $data[$col]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $row = $i + 2
        $ws.Range("$col$row").Value = $values[$i]
    }
}
